$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the existing "_GoBack" bookmark (it currently sits right
#    after "... enforce in minimizing the software debt. ").
# ------------------------------------------------------------------
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2) Merge ", " + "to ensure repeatable, high-quality, cost-effective
#    deployments of solutions" into a single run of text.
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    ", to ensure repeatable, high-quality, cost-effective deployments of solutions",
    $false, $false, $false, $false, $false, $true, 1, $false,
    ", to ensure repeatable, high-quality, cost-effective deployments of solutions",
    2) | Out-Null

# ------------------------------------------------------------------
# 3) Expand "master data management process" into
#    "master data management, Data profiling data mining process"
#    (adds the "Data mining" experience mentioned in the commit msg).
# ------------------------------------------------------------------
$d.Content.Find.Execute(
    "master data management process",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "master data management, Data profiling data mining process",
    2) | Out-Null

# ------------------------------------------------------------------
# 4) Re-insert the "_GoBack" bookmark right after the bold run
#    "integrating with Salesforce platform" (before the closing ".").
# ------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(
    "integrating with Salesforce platform",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$rng.Collapse(0)  # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $rng) | Out-Null
